$d = $word.ActiveDocument

$newText = "Stávate sa súčasťou celosvetovej kampane Globe at Night, ktorej cieľom je meranie svetelného znečistenia. Pozorovaním  Súhvezdie Blíženci na nočnej oblohe a porovnávaním skutočnej situácie s našimi mapkami sa nielenže dozviete, ako osvetlenie vo Vašom okolí prispieva k svetelnému znečisteniu, ale budete môcť porovnať úroveň svetelného znečistenia aj s inými lokalitami z celého sveta. Vaše pozorovanie tiež rozšíri online databázu dokumentujúcu viditeľnosť nočnej oblohy na našej planéte"

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Stávate sa súčasťou celosvetovej kampane*") {
        $r = $p.Range
        # Insert a throwaway marker right at the start of the paragraph's
        # content. This shifts the orphaned leading <w:proofErr/> marker so
        # that it is no longer sitting exactly on the paragraph-content
        # boundary, which lets the subsequent range delete sweep it away
        # along with the rest of the old (heavily run-split) content.
        $r0 = $d.Range($r.Start, $r.Start)
        $r0.InsertBefore("X")

        $r = $p.Range
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.Delete()
        $r2.InsertAfter($newText)
    }
}
Write-Output "DONE"
